# Certificado_Sistema_Junin.pptx -- "correcciones generales de contenido
# estadisticas y variables"
#
# This deck (as provided) contains a single content slide plus the shared
# slide master (footer / date / slide-number placeholders). The canonical
# diff nudges a batch of shape positions/sizes by a few dozen EMU (sub-pixel
# layout corrections), trims a stray leading space out of the
# "(Nombre_Comercio)" label, and tweaks a couple of paragraph/run
# properties.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Background picture -- nudge left edge in, shrink slightly.
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(1)
$sh.Left = 0.14173229038715363
$sh.Top = 0.0
$sh.Width = 719.716552734375
$sh.Height = 539.716552734375

# ---------------------------------------------------------------------
# 2) "RIF:" label rectangle.
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(2)
$sh.Left = 272.4661560058594
$sh.Top = 206.47560119628906
$sh.Width = 51.70393753051758
$sh.Height = 35.85826873779297

# ---------------------------------------------------------------------
# 3) "Direccion:" label rectangle.
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(3)
$sh.Left = 90.93544006347656
$sh.Top = 263.42364501953125
$sh.Width = 89.14961242675781
$sh.Height = 28.658267974853516

# ---------------------------------------------------------------------
# 4) "(Nombre_Comercio)" rectangle -- drop the stray leading space run
#    and reposition/resize the box.
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 1).Text = ""
$sh.Left = 260.4756164550781
$sh.Top = 165.06143188476562
$sh.Width = 187.7102508544922
$sh.Height = 31.067716598510742

# ---------------------------------------------------------------------
# 5) CuadroTexto 9 (ID_Comercio).
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(5)
$sh.Left = 616.2803344726562
$sh.Top = 425.6504211425781
$sh.Width = 143.71653747558594
$sh.Height = 21.429922103881836

# ---------------------------------------------------------------------
# 6) CuadroTexto 10 (Fecha_Emision).
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(6)
$sh.Left = 546.4630126953125
$sh.Top = 452.4945068359375
$sh.Width = 98.27716827392578
$sh.Height = 28.658267974853516

# ---------------------------------------------------------------------
# 7) CuadroTexto 12 (Rif_Empresarial).
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(7)
$sh.Left = 313.71026611328125
$sh.Top = 208.0063018798828
$sh.Width = 139.66299438476562
$sh.Height = 28.658267974853516

# ---------------------------------------------------------------------
# 8) CuadroTexto 14 (Direccion).
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(8)
$sh.Left = 168.40631103515625
$sh.Top = 263.6504211425781
$sh.Width = 384.4630126953125
$sh.Height = 28.658267974853516

# ---------------------------------------------------------------------
# 9) QR placeholder square behind the "[Q]" box.
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(9)
$sh.Left = 581.1024169921875
$sh.Top = 311.81103515625
$sh.Width = 113.21575164794922
$sh.Height = 113.21575164794922
$sh.TextFrame.TextRange.ParagraphFormat.SpaceWithin = 1

# ---------------------------------------------------------------------
# 10) "[Q]" rectangle.
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(10)
$sh.Left = 575.5181274414062
$sh.Top = 286.2992248535156
$sh.Width = 121.63465118408203
$sh.Height = 141.5622100830078

# ---------------------------------------------------------------------
# 11) Slide master -- footer / slide-number / date placeholders shrink
#     by a couple EMU to match the corrected layout.
# ---------------------------------------------------------------------
$m = $p.SlideMaster

$mf = $m.Shapes.Item(1)
$mf.Left = 238.50709533691406
$mf.Top = 500.5133972167969
$mf.Width = 242.730712890625
$mf.Height = 28.459842681884766

$mn = $m.Shapes.Item(2)
$mn.Left = 508.5071105957031
$mn.Top = 500.5133972167969
$mn.Width = 161.71653747558594
$mn.Height = 28.459842681884766

$md = $m.Shapes.Item(3)
$md.Left = 49.492916107177734
$md.Top = 500.5133972167969
$md.Width = 161.71653747558594
$md.Height = 28.459842681884766
